{"js": "const pairs = [\n  [\"6+75=81\", \"4+77=81\"],\n  [\"44-12=32\", \"67+22=89\"],\n  [\"31+11=42\", \"75+1=76\"],\n  [\"64-49=15\", \"21+68=89\"],\n  [\"20+44=64\", \"55-10=45\"],\n  [\"41+42=83\", \"14+23=37\"],\n  [\"94-84=10\", \"68-8=60\"],\n  [\"71+8=79\", \"93-83=10\"],\n  [\"30-12=18\", \"7+89=96\"],\n  [\"51-42=9\", \"27+6=33\"],\n  [\"34-2=32\", \"7+21=28\"],\n  [\"29+12=41\", \"27-13=14\"],\n  [\"29+13=42\", \"21+74=95\"],\n  [\"98-38=60\", \"52+2=54\"],\n  [\"1+20=21\", \"39+10=49\"],\n  [\"73+1=74\", \"53-29=24\"],\n  [\"70-65=5\", \"70+17=87\"],\n  [\"86-51=35\", \"22+40=62\"],\n  [\"61-45=16\", \"98-30=68\"],\n  [\"34+39=73\", \"18+0=18\"],\n  [\"43+56=99\", \"47-41=6\"],\n  [\"74-56=18\", \"35+60=95\"],\n  [\"32-7=25\", \"70+24=94\"],\n  [\"9+72=81\", \"46+17=63\"],\n  [\"26+68=94\", \"71-1=70\"],\n  [\"57-15=42\", \"47+21=68\"],\n  [\"69-3=66\", \"52-17=35\"],\n  [\"82-0=82\", \"87-49=38\"],\n  [\"4+54=58\", \"48-47=1\"],\n  [\"93+5=98\", \"73-47=26\"],\n  [\"6-0=6\", \"11+61=72\"],\n  [\"56-18=38\", \"53-43=10\"],\n  [\"65-41=24\", \"55-39=16\"],\n  [\"39-4=35\", \"88+11=99\"],\n  [\"85-9=76\", \"2+87=89\"],\n  [\"20+73=93\", \"13-12=1\"],\n  [\"75+15=90\", \"87-75=12\"],\n  [\"45-34=11\", \"80-39=41\"],\n  [\"3+86=89\", \"71-13=58\"],\n  [\"6-3=3\", \"15+53=68\"],\n  [\"47+7=54\", \"23+33=56\"],\n  [\"67-22=45\", \"86-30=56\"],\n  [\"61-5=56\", \"66+2=68\"],\n  [\"69+24=93\", \"0+6=6\"],\n  [\"94-38=56\", \"60-33=27\"],\n  [\"51+14=65\", \"53+30=83\"],\n  [\"57+21=78\", \"45+47=92\"],\n  [\"37+32=69\", \"46-14=32\"],\n  [\"51-30=21\", \"90-24=66\"],\n  [\"25+27=52\", \"76-14=62\"],\n  [\"81+14=95\", \"29+11=40\"],\n  [\"23+53=76\", \"70-8=62\"],\n  [\"67-14=53\", \"3+61=64\"],\n  [\"50+14=64\", \"58+4=62\"],\n  [\"42-33=9\", \"66+33=99\"],\n  [\"40-39=1\", \"73+7=80\"],\n  [\"70-66=4\", \"47-39=8\"],\n  [\"8+58=66\", \"58-29=29\"],\n  [\"94-50=44\", \"79-9=70\"],\n  [\"7+46=53\", \"21-17=4\"],\n  [\"55-31=24\", \"41-38=3\"],\n  [\"50+15=65\", \"76-41=35\"],\n  [\"75-1=74\", \"22+23=45\"],\n  [\"61-42=19\", \"69-62=7\"],\n  [\"5+42=47\", \"42-15=27\"],\n  [\"33+4=37\", \"55-1=54\"],\n  [\"82-9=73\", \"18+55=73\"],\n  [\"76-9=67\", \"65-11=54\"],\n  [\"40+22=62\", \"9+87=96\"],\n  [\"30+41=71\", \"74+12=86\"],\n  [\"99-25=74\", \"59-39=20\"],\n  [\"6+88=94\", \"89-14=75\"],\n  [\"8+48=56\", \"31+43=74\"],\n  [\"35+11=46\", \"69+1=70\"],\n  [\"14+77=91\", \"2+81=83\"],\n  [\"56+18=74\", \"95-50=45\"],\n  [\"85-1=84\", \"10+5=15\"],\n  [\"99-75=24\", \"6+24=30\"],\n  [\"52-52=0\", \"51+30=81\"],\n  [\"17-15=2\", \"90-12=78\"],\n  [\"87-76=11\", \"51-24=27\"],\n  [\"99-95=4\", \"63+10=73\"],\n  [\"75-22=53\", \"37+61=98\"],\n  [\"26+51=77\", \"20+52=72\"],\n  [\"27-9=18\", \"67-21=46\"],\n  [\"78-71=7\", \"67+14=81\"],\n  [\"32+10=42\", \"43+42=85\"],\n  [\"11+1=12\", \"16-16=0\"],\n  [\"94-29=65\", \"30+5=35\"],\n  [\"44-30=14\", \"37+38=75\"],\n  [\"97-8=89\", \"92-13=79\"],\n  [\"45-45=0\", \"98-29=69\"],\n  [\"5+87=92\", \"33+28=61\"],\n  [\"45+54=99\", \"2+61=63\"],\n  [\"92-75=17\", \"31-29=2\"],\n  [\"48-41=7\", \"11+79=90\"],\n  [\"33+35=68\", \"39+60=99\"],\n  [\"64-13=51\", \"83-77=6\"],\n  [\"99-39=60\", \"94-34=60\"],\n  [\"45-30=15\", \"64-52=12\"],\n];\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const result of results.items) {\n    result.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$pairs = @(\n  @(\"6+75=81\", \"4+77=81\"),\n  @(\"44-12=32\", \"67+22=89\"),\n  @(\"31+11=42\", \"75+1=76\"),\n  @(\"64-49=15\", \"21+68=89\"),\n  @(\"20+44=64\", \"55-10=45\"),\n  @(\"41+42=83\", \"14+23=37\"),\n  @(\"94-84=10\", \"68-8=60\"),\n  @(\"71+8=79\", \"93-83=10\"),\n  @(\"30-12=18\", \"7+89=96\"),\n  @(\"51-42=9\", \"27+6=33\"),\n  @(\"34-2=32\", \"7+21=28\"),\n  @(\"29+12=41\", \"27-13=14\"),\n  @(\"29+13=42\", \"21+74=95\"),\n  @(\"98-38=60\", \"52+2=54\"),\n  @(\"1+20=21\", \"39+10=49\"),\n  @(\"73+1=74\", \"53-29=24\"),\n  @(\"70-65=5\", \"70+17=87\"),\n  @(\"86-51=35\", \"22+40=62\"),\n  @(\"61-45=16\", \"98-30=68\"),\n  @(\"34+39=73\", \"18+0=18\"),\n  @(\"43+56=99\", \"47-41=6\"),\n  @(\"74-56=18\", \"35+60=95\"),\n  @(\"32-7=25\", \"70+24=94\"),\n  @(\"9+72=81\", \"46+17=63\"),\n  @(\"26+68=94\", \"71-1=70\"),\n  @(\"57-15=42\", \"47+21=68\"),\n  @(\"69-3=66\", \"52-17=35\"),\n  @(\"82-0=82\", \"87-49=38\"),\n  @(\"4+54=58\", \"48-47=1\"),\n  @(\"93+5=98\", \"73-47=26\"),\n  @(\"6-0=6\", \"11+61=72\"),\n  @(\"56-18=38\", \"53-43=10\"),\n  @(\"65-41=24\", \"55-39=16\"),\n  @(\"39-4=35\", \"88+11=99\"),\n  @(\"85-9=76\", \"2+87=89\"),\n  @(\"20+73=93\", \"13-12=1\"),\n  @(\"75+15=90\", \"87-75=12\"),\n  @(\"45-34=11\", \"80-39=41\"),\n  @(\"3+86=89\", \"71-13=58\"),\n  @(\"6-3=3\", \"15+53=68\"),\n  @(\"47+7=54\", \"23+33=56\"),\n  @(\"67-22=45\", \"86-30=56\"),\n  @(\"61-5=56\", \"66+2=68\"),\n  @(\"69+24=93\", \"0+6=6\"),\n  @(\"94-38=56\", \"60-33=27\"),\n  @(\"51+14=65\", \"53+30=83\"),\n  @(\"57+21=78\", \"45+47=92\"),\n  @(\"37+32=69\", \"46-14=32\"),\n  @(\"51-30=21\", \"90-24=66\"),\n  @(\"25+27=52\", \"76-14=62\"),\n  @(\"81+14=95\", \"29+11=40\"),\n  @(\"23+53=76\", \"70-8=62\"),\n  @(\"67-14=53\", \"3+61=64\"),\n  @(\"50+14=64\", \"58+4=62\"),\n  @(\"42-33=9\", \"66+33=99\"),\n  @(\"40-39=1\", \"73+7=80\"),\n  @(\"70-66=4\", \"47-39=8\"),\n  @(\"8+58=66\", \"58-29=29\"),\n  @(\"94-50=44\", \"79-9=70\"),\n  @(\"7+46=53\", \"21-17=4\"),\n  @(\"55-31=24\", \"41-38=3\"),\n  @(\"50+15=65\", \"76-41=35\"),\n  @(\"75-1=74\", \"22+23=45\"),\n  @(\"61-42=19\", \"69-62=7\"),\n  @(\"5+42=47\", \"42-15=27\"),\n  @(\"33+4=37\", \"55-1=54\"),\n  @(\"82-9=73\", \"18+55=73\"),\n  @(\"76-9=67\", \"65-11=54\"),\n  @(\"40+22=62\", \"9+87=96\"),\n  @(\"30+41=71\", \"74+12=86\"),\n  @(\"99-25=74\", \"59-39=20\"),\n  @(\"6+88=94\", \"89-14=75\"),\n  @(\"8+48=56\", \"31+43=74\"),\n  @(\"35+11=46\", \"69+1=70\"),\n  @(\"14+77=91\", \"2+81=83\"),\n  @(\"56+18=74\", \"95-50=45\"),\n  @(\"85-1=84\", \"10+5=15\"),\n  @(\"99-75=24\", \"6+24=30\"),\n  @(\"52-52=0\", \"51+30=81\"),\n  @(\"17-15=2\", \"90-12=78\"),\n  @(\"87-76=11\", \"51-24=27\"),\n  @(\"99-95=4\", \"63+10=73\"),\n  @(\"75-22=53\", \"37+61=98\"),\n  @(\"26+51=77\", \"20+52=72\"),\n  @(\"27-9=18\", \"67-21=46\"),\n  @(\"78-71=7\", \"67+14=81\"),\n  @(\"32+10=42\", \"43+42=85\"),\n  @(\"11+1=12\", \"16-16=0\"),\n  @(\"94-29=65\", \"30+5=35\"),\n  @(\"44-30=14\", \"37+38=75\"),\n  @(\"97-8=89\", \"92-13=79\"),\n  @(\"45-45=0\", \"98-29=69\"),\n  @(\"5+87=92\", \"33+28=61\"),\n  @(\"45+54=99\", \"2+61=63\"),\n  @(\"92-75=17\", \"31-29=2\"),\n  @(\"48-41=7\", \"11+79=90\"),\n  @(\"33+35=68\", \"39+60=99\"),\n  @(\"64-13=51\", \"83-77=6\"),\n  @(\"99-39=60\", \"94-34=60\"),\n  @(\"45-30=15\", \"64-52=12\"),\n)\n$notFound = @()\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $range = $d.Content\n  $find = $range.Find\n  $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    $notFound += $oldText\n  }\n}\n\nif ($notFound.Count -gt 0) {\n  \"NOT FOUND: \" + ($notFound -join \", \")\n} else {\n  \"All replacements applied successfully\"\n}\n"}
